$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A12").Value = "BASESECRETURI"
$ws.Range("B12").Value = "https://cicdsecretkeyvalues.vault.azure.net/"
$ws.Range("C12").Value = "DNS value of Azure Secret key vault"
$ws.Range("C13").Select() | Out-Null
